$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws1 = $wb.Worksheets.Item("Rushing")

# Row 2 - M.Ryan
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 8
$ws1.Range("E2").Value = 12
$ws1.Range("F2").Value = 5

# Row 4 - M.Davis
$ws1.Range("C4").Value = 79
$ws1.Range("D4").Value = 30
$ws1.Range("E4").Value = 8
$ws1.Range("F4").Value = 14

# Row 5 - C.Patterson
$ws1.Range("C5").Value = 90
$ws1.Range("D5").Value = 34
$ws1.Range("E5").Value = 10
$ws1.Range("F5").Value = 28

# Row 7 - Q.Ollison
$ws1.Range("C7").Value = 13

# --- Receiving sheet ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 2 - M.Davis
$ws2.Range("C2").Value = 50
$ws2.Range("G2").Value = 6

# Row 3 - C.Patterson
$ws2.Range("C3").Value = 54
$ws2.Range("D3").Value = 40
$ws2.Range("G3").Value = 12
$ws2.Range("H3").Value = 9

# Row 5 - Q.Ollison
$ws2.Range("C5").Value = 5
$ws2.Range("D5").Value = 4
$ws2.Range("G5").Value = 1
$ws2.Range("H5").Value = 1

# Row 6 - R.Gage
$ws2.Range("C6").Value = 69
$ws2.Range("D6").Value = 58
$ws2.Range("E6").Value = 13
$ws2.Range("F6").Value = 9
$ws2.Range("G6").Value = 9
$ws2.Range("H6").Value = 6

# Row 7 - O.Zaccheaus
$ws2.Range("C7").Value = 35
$ws2.Range("D7").Value = 20
$ws2.Range("E7").Value = 6
$ws2.Range("F7").Value = 3
$ws2.Range("G7").Value = 7
$ws2.Range("H7").Value = 4

# Row 8 - C.Blake
$ws2.Range("C8").Value = 5
$ws2.Range("G8").Value = 1

# Row 9 - T.Sharpe
$ws2.Range("C9").Value = 29
$ws2.Range("G9").Value = 3

# Row 10 - K.Pitts
$ws2.Range("C10").Value = 69
$ws2.Range("D10").Value = 44
$ws2.Range("E10").Value = 25
$ws2.Range("F10").Value = 14
$ws2.Range("G10").Value = 13
$ws2.Range("H10").Value = 4

# Row 11 - H.Hurst
$ws2.Range("C11").Value = 27
$ws2.Range("D11").Value = 23

# Row 12 - L.Smith
$ws2.Range("C12").Value = 10
$ws2.Range("G12").Value = 2
